$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Periodo Mora" data column (E16:E19) so the period codes are
# listed in ascending order (2012, 2101, 2102, 2103) instead of the prior
# descending order (2103, 2102, 2101, 2012). Cells are text-formatted
# ("@" number format), so the values are written as strings.
$ws.Range("E16").Value = "2012"
$ws.Range("E17").Value = "2101"
$ws.Range("E18").Value = "2102"
$ws.Range("E19").Value = "2103"
